$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert a new column before the existing "sliderStep" column (W), which
#    shifts sliderStep/paramHTML/multi-slider?/sliderFun/transformFunList/
#    pdfList/distrPlotList/randomDrawList/latexList/chartDomain/likelihoodFun
#    one column to the right (W->X, X->Y, ... AG->AH). Excel auto-adjusts all
#    formula references (relative + $-anchored) that point past the new
#    column, and the new column inherits the "Neutral" style from its left
#    neighbour (sigmaScale, column V).
# ---------------------------------------------------------------------------
$ws.Columns("W").Insert()

# ---------------------------------------------------------------------------
# 2) Insert a new row before row 7 (the old "styNorm" row), pushing styNorm
#    and everything below it down by one row (old 7-20 -> new 8-21). The new
#    row will hold the "Ordered Logit (X)" distribution, a sibling of the
#    "Ordered Probit (X)" row directly above it (row 6).
# ---------------------------------------------------------------------------
$ws.Rows(7).Insert()

# ---------------------------------------------------------------------------
# 3) Header for the newly inserted column W: "yStarPDF"
# ---------------------------------------------------------------------------
$ws.Range("W1").Value = "yStarPDF"

# ---------------------------------------------------------------------------
# 4) Fill the new column W ("yStarPDF") for every data row. It is "NA" for
#    every distribution except the two ordered-outcome rows, which reference
#    the shared y*-PDF helper used by their respective link functions.
# ---------------------------------------------------------------------------
$ws.Range("W2").Value = "NA"
$ws.Range("W3").Value = "NA"
$ws.Range("W4").Value = "NA"
$ws.Range("W5").Value = "NA"
$ws.Range("W6").Value = "styNormPDF"
$ws.Range("W7").Value = "styLogPDF"
$ws.Range("W8").Value = "NA"
$ws.Range("W9").Value = "NA"
$ws.Range("W10").Value = "NA"
$ws.Range("W11").Value = "NA"
$ws.Range("W12").Value = "NA"
$ws.Range("W13").Value = "NA"
$ws.Range("W14").Value = "NA"
$ws.Range("W15").Value = "NA"
$ws.Range("W16").Value = "NA"
$ws.Range("W17").Value = "NA"
$ws.Range("W18").Value = "NA"
$ws.Range("W19").Value = "NA"
$ws.Range("W20").Value = "NA"
$ws.Range("W21").Value = "NA"

# ---------------------------------------------------------------------------
# 5) Row 6 ("orderedProbitX") now shares its distrGroups with the new row,
#    so its group label changes from "Ordered Probit (X)" to "Ordered".
# ---------------------------------------------------------------------------
$ws.Range("D6").Value = "Ordered"

# ---------------------------------------------------------------------------
# 6) Populate the new row 7 ("orderedLogitX") - first draft of the ordered
#    logit model, mirroring row 6 (orderedProbitX) with logit-specific
#    names/links.
# ---------------------------------------------------------------------------
$ws.Range("A7").Formula = "=A6+1"
$ws.Range("B7").Value = "orderedLogitX"
$ws.Range("C7").Value = "Ordered Logit (X)"
$ws.Range("D7").Value = "Ordered"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 2
$ws.Range("G7").Value = "intPrintHelper"
$ws.Range("H7").Value = 'list("Ordered Logit (X)")'
$ws.Range("I7").Value = 'c("Beta0", "Beta1", "Beta2")'
$ws.Range("J7").Value = "L-BFGS-B"
$ws.Range("K7").Value = '$ \tilde{E}(y) =\tilde{\pi} = \tilde{Pr}(Y=1)$'
$ws.Range("L7").Value = "Beta"
$ws.Range("M7").Value = '\beta'
$ws.Range("N7").Value = '\pi'
$ws.Range("O7").Value = 'list("Predicted Values", "Expected Values")'
$ws.Range("P7").Value = "c(1,3)"
$ws.Range("Q7").Value = "NA"
$ws.Range("R7").Value = "c(0,1)"
$ws.Range("S7").Value = -3
$ws.Range("T7").Value = 3
$ws.Range("U7").Value = "c(1,-1,1.25)"
$ws.Range("V7").Value = "c(-1.5,1.5)"
$ws.Range("W7").Value = "styLogPDF"
$ws.Range("Y7").Formula = '="""&"&RIGHT(M7,LEN(M7)-1)&";"""'
$ws.Range("Z7").Formula = '=IF(F7=1,"""none""",IF(E7=F7,"""betas""","""fullNorm"""))'
$ws.Range("AA7").Formula = '="manyParamSliderMaker(minVal ="&S7&", maxVal = "&T7&", startVals = "&U7&", stepVal = "&X7&", paramHTML = "&Y7&", multi = "&Z7&", sigmaScale ="&V7&","'
$ws.Range("AB7").Formula = '=$B7&"ParamTransform"'
$ws.Range("AC7").Formula = '=$B7&"PDF"'
$ws.Range("AD7").Formula = '=B7&"PlotDistr"'
$ws.Range("AE7").Formula = '=B7&"Draws"'
$ws.Range("AF7").Formula = '=B7&"Latex"'
$ws.Range("AG7").Formula = '=$B7&"ChartDomain"'
$ws.Range("AH7").Formula = '=$B7&"LikelihoodFun"'

# ---------------------------------------------------------------------------
# 7) Restore the selection to match the authored state.
# ---------------------------------------------------------------------------
$ws.Range("W8").Select()
